$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("set14")

# Update the lunch-break times on row 8 (D8 = almoço saída, G8 = almoço entrada)
$ws.Range("D8").Value = 0.53333333333333333
$ws.Range("G8").Value = 0.79375000000000007

# Update the selected cell shown when the workbook is opened
$ws.Range("Q7").Select()
